$wb = $excel.ActiveWorkbook

# --- Rename sheets ------------------------------------------------------
$wsMse1 = $wb.Worksheets.Item("Sheet1")
$wsMse1.Name = "MSE_1"

$ws = $wb.Worksheets.Item("MSE")
$ws.Name = "MSE_2"

# --- Re-layout the MSE_2 sheet ------------------------------------------
# Push the existing TV/radio/newspaper/sales/f_hat1/f_hat2 data (and the
# page-setup / styling that rides along with it) six columns to the right,
# out of the way, so the cross-validation predictions can take over the
# front A:C columns.
$ws.Columns("A:F").Insert() | Out-Null

# The y_true / y_pred1 / y_pred2 columns (now shifted to R:T by the insert
# above) become the new headline columns A:C.
$ws.Range("R1:T6").Cut($ws.Range("A1")) | Out-Null

# --- Hide the now-secondary feature/helper columns -----------------------
$ws.Columns("F:M").ColumnWidth = 0
$ws.Columns("F:M").Hidden = $true

# --- Selection / active sheet --------------------------------------------
$ws.Range("D1:E1048576").Select() | Out-Null
